# Login details updated in Test data sheet and share skills edit details
# entered in share skill excel sheet.
#
# The SignIn ("Test data") sheet's sample login row is refreshed with new
# credentials. (The SignUp / Profile sheets keep their existing values —
# only their underlying shared-string indices shift because the old,
# now-unused URL string gets dropped from the shared string table.)

$wb = $excel.ActiveWorkbook

$signIn = $wb.Worksheets.Item("SignIn")

# New Url / Username / Password for the SignIn test-data row.
# Write Password (C2) first, then Username (B2), then Url (A2) so the new
# shared-string entries are appended in that same order.
$signIn.Range("C2").Value = "Mars@123"
$signIn.Range("B2").Value = "archika.mehta19@gmail.com"
$signIn.Range("A2").Value = "http://192.168.99.100:5000"

# Make SignIn the active/selected sheet with A2 highlighted, matching the
# saved workbook view state.
$signIn.Activate()
$signIn.Range("A2").Select()
